# Updated cryptos list values (Price / Volume(1h)) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.028.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.11%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.869.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.11%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.40%  "

$ws.Range("E6").Value = "  +0.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5101"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.59%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3801"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.36%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08301"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -9.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.110"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.67%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.44"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.78%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.212"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.49%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.869.05"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.13%  "

$ws.Range("E14").Value = "  -1.73%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.179"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.17%  "

$ws.Range("E17").Value = "  -1.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.68"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06628"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.05%  "

$ws.Range("E20").Value = "  -0.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.60%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.062.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.258"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.72%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.568"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.57%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.084.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.44%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "157.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.84%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.46%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.18%  "

$ws.Range("E31").Value = "  +0.21%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.040"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.18%  "

$ws.Range("E33").Value = "  -0.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.607"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.669"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.86%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02430"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.67%  "

$ws.Range("E37").Value = "  -0.69%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2160"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.66%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.206"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.87%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6419"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.39%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.239"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.74%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.74%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.870"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.65%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6100"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.70%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.284"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.662"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.67%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.996"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.18%  "

$ws.Range("E49").Value = "  +1.23%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "121.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.43%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.86"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.99%  "
